$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule "R30" (row 10): update the "Integer min" (From) value in column C
# from 18 to 1.
$ws.Range("C10").Value = 1
